# Edit script: updates the date heading and every answer cell in the
# single 20x5 table to match the target revision.
$d = $word.ActiveDocument

# --- Update the date heading (first paragraph) ---
$d.Paragraphs.Item(1).Range.Text = "2024-01-31 Wednesday"

# --- Update each cell of the single table, in row-major order ---
# Cell.Range.Text is set directly (positional), not Find&Replace, because
# several old values are not unique within the table and map to different
# new values depending on which occurrence they are.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "19+43=62"  # was: 23+69=92
$tbl.Cell(1, 2).Range.Text = "40-6=34"  # was: 49+3=52
$tbl.Cell(1, 3).Range.Text = "9+74=83"  # was: 43-37=6
$tbl.Cell(1, 4).Range.Text = "84-56=28"  # was: 19+43=62
$tbl.Cell(1, 5).Range.Text = "41-3=38"  # was: 19+64=83

$tbl.Cell(2, 1).Range.Text = "61-37=24"  # was: 93-89=4
$tbl.Cell(2, 2).Range.Text = "67+26=93"  # was: 46-19=27
$tbl.Cell(2, 3).Range.Text = "41-18=23"  # was: 18+26=44
$tbl.Cell(2, 4).Range.Text = "9+86=95"  # was: 59+37=96
$tbl.Cell(2, 5).Range.Text = "18+59=77"  # was: 77-69=8

$tbl.Cell(3, 1).Range.Text = "86+5=91"  # was: 38+53=91
$tbl.Cell(3, 2).Range.Text = "30-2=28"  # was: 60-11=49
$tbl.Cell(3, 3).Range.Text = "80-53=27"  # was: 19+63=82
$tbl.Cell(3, 4).Range.Text = "26+29=55"  # was: 16+8=24
$tbl.Cell(3, 5).Range.Text = "82-37=45"  # was: 56-7=49

$tbl.Cell(4, 1).Range.Text = "39+19=58"  # was: 50-12=38
$tbl.Cell(4, 2).Range.Text = "18+43=61"  # was: 96-18=78
$tbl.Cell(4, 3).Range.Text = "76+8=84"  # was: 8+53=61
$tbl.Cell(4, 4).Range.Text = "31-2=29"  # was: 57+28=85
$tbl.Cell(4, 5).Range.Text = "61-24=37"  # was: 48+17=65

$tbl.Cell(5, 1).Range.Text = "80-36=44"  # was: 70-4=66
$tbl.Cell(5, 2).Range.Text = "51-2=49"  # was: 40-17=23
$tbl.Cell(5, 3).Range.Text = "60-29=31"  # was: 35+16=51
$tbl.Cell(5, 4).Range.Text = "63+8=71"  # was: 86-17=69
$tbl.Cell(5, 5).Range.Text = "49+46=95"  # was: 55+27=82

$tbl.Cell(6, 1).Range.Text = "59+5=64"  # was: 93-58=35
$tbl.Cell(6, 2).Range.Text = "61-54=7"  # was: 58+28=86
$tbl.Cell(6, 3).Range.Text = "39+46=85"  # was: 31-22=9
$tbl.Cell(6, 4).Range.Text = "34+39=73"  # was: 19+76=95
$tbl.Cell(6, 5).Range.Text = "50-33=17"  # was: 70-46=24

$tbl.Cell(7, 1).Range.Text = "31-27=4"  # was: 74-58=16
$tbl.Cell(7, 2).Range.Text = "90-68=22"  # was: 29+65=94
$tbl.Cell(7, 3).Range.Text = "24+9=33"  # was: 4+58=62
$tbl.Cell(7, 4).Range.Text = "6+36=42"  # was: 90-44=46
$tbl.Cell(7, 5).Range.Text = "82-18=64"  # was: 22+49=71

$tbl.Cell(8, 1).Range.Text = "32-28=4"  # was: 73-34=39
$tbl.Cell(8, 2).Range.Text = "81-62=19"  # was: 58+7=65
$tbl.Cell(8, 3).Range.Text = "39+45=84"  # was: 56+39=95
$tbl.Cell(8, 4).Range.Text = "91-47=44"  # was: 61-43=18
$tbl.Cell(8, 5).Range.Text = "38+8=46"  # was: 93-55=38

$tbl.Cell(9, 1).Range.Text = "23-15=8"  # was: 17-8=9
$tbl.Cell(9, 2).Range.Text = "39+45=84"  # was: 98-89=9
$tbl.Cell(9, 3).Range.Text = "82-79=3"  # was: 27+9=36
$tbl.Cell(9, 4).Range.Text = "9+59=68"  # was: 84-69=15
$tbl.Cell(9, 5).Range.Text = "17+57=74"  # was: 92-77=15

$tbl.Cell(10, 1).Range.Text = "29+13=42"  # was: 45+6=51
$tbl.Cell(10, 2).Range.Text = "40-37=3"  # was: 91-79=12
$tbl.Cell(10, 3).Range.Text = "17+57=74"  # was: 31-22=9
$tbl.Cell(10, 4).Range.Text = "46+15=61"  # was: 27+36=63
$tbl.Cell(10, 5).Range.Text = "66+26=92"  # was: 91-19=72

$tbl.Cell(11, 1).Range.Text = "76+8=84"  # was: 16+78=94
$tbl.Cell(11, 2).Range.Text = "43+48=91"  # was: 36-29=7
$tbl.Cell(11, 3).Range.Text = "41-16=25"  # was: 30-21=9
$tbl.Cell(11, 4).Range.Text = "91-14=77"  # was: 21-19=2
$tbl.Cell(11, 5).Range.Text = "29+25=54"  # was: 49+12=61

$tbl.Cell(12, 1).Range.Text = "16+46=62"  # was: 29+57=86
$tbl.Cell(12, 2).Range.Text = "65-56=9"  # was: 18+13=31
$tbl.Cell(12, 3).Range.Text = "82+9=91"  # was: 69+9=78
$tbl.Cell(12, 4).Range.Text = "72-26=46"  # was: 16+28=44
$tbl.Cell(12, 5).Range.Text = "75-67=8"  # was: 60-34=26

$tbl.Cell(13, 1).Range.Text = "37+28=65"  # was: 27+19=46
$tbl.Cell(13, 2).Range.Text = "38+9=47"  # was: 97-58=39
$tbl.Cell(13, 3).Range.Text = "86-19=67"  # was: 9+28=37
$tbl.Cell(13, 4).Range.Text = "32-13=19"  # was: 81-17=64
$tbl.Cell(13, 5).Range.Text = "57+27=84"  # was: 17+69=86

$tbl.Cell(14, 1).Range.Text = "22-6=16"  # was: 74+18=92
$tbl.Cell(14, 2).Range.Text = "18+37=55"  # was: 4+27=31
$tbl.Cell(14, 3).Range.Text = "85-38=47"  # was: 81-29=52
$tbl.Cell(14, 4).Range.Text = "58-9=49"  # was: 8+35=43
$tbl.Cell(14, 5).Range.Text = "48+26=74"  # was: 48+23=71

$tbl.Cell(15, 1).Range.Text = "83-68=15"  # was: 27+27=54
$tbl.Cell(15, 2).Range.Text = "19+34=53"  # was: 81-78=3
$tbl.Cell(15, 3).Range.Text = "15+58=73"  # was: 57+24=81
$tbl.Cell(15, 4).Range.Text = "70-2=68"  # was: 40-22=18
$tbl.Cell(15, 5).Range.Text = "67+26=93"  # was: 63-55=8

$tbl.Cell(16, 1).Range.Text = "83-54=29"  # was: 64-45=19
$tbl.Cell(16, 2).Range.Text = "6+55=61"  # was: 46+9=55
$tbl.Cell(16, 3).Range.Text = "66+7=73"  # was: 62-16=46
$tbl.Cell(16, 4).Range.Text = "91-27=64"  # was: 39+26=65
$tbl.Cell(16, 5).Range.Text = "80-2=78"  # was: 33+38=71

$tbl.Cell(17, 1).Range.Text = "7+58=65"  # was: 64-36=28
$tbl.Cell(17, 2).Range.Text = "11-7=4"  # was: 91-25=66
$tbl.Cell(17, 3).Range.Text = "10-2=8"  # was: 9+48=57
$tbl.Cell(17, 4).Range.Text = "90-58=32"  # was: 24-6=18
$tbl.Cell(17, 5).Range.Text = "64-56=8"  # was: 56-49=7

$tbl.Cell(18, 1).Range.Text = "52-35=17"  # was: 92-8=84
$tbl.Cell(18, 2).Range.Text = "91-45=46"  # was: 83-9=74
$tbl.Cell(18, 3).Range.Text = "90-32=58"  # was: 70-55=15
$tbl.Cell(18, 4).Range.Text = "25+9=34"  # was: 61-42=19
$tbl.Cell(18, 5).Range.Text = "19+43=62"  # was: 19+12=31

$tbl.Cell(19, 1).Range.Text = "36+36=72"  # was: 15+18=33
$tbl.Cell(19, 2).Range.Text = "87-28=59"  # was: 18+54=72
$tbl.Cell(19, 3).Range.Text = "87-39=48"  # was: 20-4=16
$tbl.Cell(19, 4).Range.Text = "18+23=41"  # was: 28+66=94
$tbl.Cell(19, 5).Range.Text = "19+4=23"  # was: 93-55=38

$tbl.Cell(20, 1).Range.Text = "55-27=28"  # was: 58+5=63
$tbl.Cell(20, 2).Range.Text = "62-8=54"  # was: 85+8=93
$tbl.Cell(20, 3).Range.Text = "45+9=54"  # was: 68+29=97
$tbl.Cell(20, 4).Range.Text = "3+38=41"  # was: 98-59=39
$tbl.Cell(20, 5).Range.Text = "34+7=41"  # was: 9+58=67

Write-Output "Done."